$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new taxonomy header columns (D:J) ---
$ws.Range("D1").Value = "Kingdom"
$ws.Range("E1").Value = "Phylum"
$ws.Range("F1").Value = "Class"
$ws.Range("G1").Value = "Order"
$ws.Range("H1").Value = "Family"
$ws.Range("I1").Value = "Genus"
$ws.Range("J1").Value = "species"

# Match the header formatting (bold, centered) already used for A1:C1
$ws.Range("A1:C1").Copy()
$ws.Range("D1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2: Fundulus heteroclitus or majalis (ambiguous -> no single species epithet) ---
$ws.Range("D2").Value = "Animalia"
$ws.Range("E2").Value = "Chordata"
$ws.Range("F2").Value = "Teleostei"
$ws.Range("G2").Value = "Cyprinodontiformes"
$ws.Range("H2").Value = "Fundulidae"
$ws.Range("I2").Value = "Fundulus"

# --- Row 3: Cololabis saira ---
$ws.Range("D3").Value = "Animalia"
$ws.Range("E3").Value = "Chordata"
$ws.Range("F3").Value = "Teleostei"
$ws.Range("G3").Value = "Beloniformes"
$ws.Range("H3").Value = "Scomberesocidae"
$ws.Range("I3").Value = "Cololabis"
$ws.Range("J3").Value = "saira"

# --- Row 4: Unassigned ---
$ws.Range("D4").Value = "Unassigned"
$ws.Range("E4").Value = "Unassigned"
$ws.Range("F4").Value = "Unassigned"
$ws.Range("G4").Value = "Unassigned"
$ws.Range("H4").Value = "Unassigned"
$ws.Range("I4").Value = "Unassigned"
$ws.Range("J4").Value = "Unassigned"

# --- Row 5: Mareca americana ---
$ws.Range("D5").Value = "Animalia"
$ws.Range("E5").Value = "Chordata"
$ws.Range("F5").Value = "Aves"
$ws.Range("G5").Value = "Anseriformes"
$ws.Range("H5").Value = "Anatidae"
$ws.Range("I5").Value = "Mareca"
$ws.Range("J5").Value = "americana"

# --- Row 6: Myrophis vafer ---
$ws.Range("D6").Value = "Animalia"
$ws.Range("E6").Value = "Chordata"
$ws.Range("F6").Value = "Teleostei"
$ws.Range("G6").Value = "Anguilliformes"
$ws.Range("H6").Value = "Ophichthidae"
$ws.Range("I6").Value = "Myrophis"
$ws.Range("J6").Value = "vafer"

# --- Column widths: narrow column A to fit the species-name text, and set a
#     uniform, narrower width for the new taxonomy columns B:J ---
$ws.Columns.Item(1).ColumnWidth = 26.2
$ws.Range("B1:J1").EntireColumn.ColumnWidth = 10.7

# --- Selection, mirroring where the user last clicked in the sheet ---
$ws.Range("I14").Select()
